$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Remove hyperlinks first (rows 8-10 had hyperlinked cells in F/G)
$ws.Hyperlinks.Delete()

# Delete the now-obsolete rows 9 and 10 (Presupuestal / Programático entries)
$ws.Rows("9:10").Delete()

# Update row 8: new reporting period + validation/update dates
$ws.Range("B8").Value = 44652
$ws.Range("C8").Value = 44742
$ws.Range("I8").Value = 44753
$ws.Range("J8").Value = 44753

# Clear the document-type / description / links columns for row 8
$ws.Range("D8:G8").ClearContents()

# Area responsible changes to Subdireccion de Planeacion y Presupuesto (UPP)
$ws.Range("H8").Value = "Subdirección de Planeación y Presupuesto (UPP)"

# New note explaining the late disclosure
$ws.Range("K8").Value = "La información solicitada se entregará a más tardar el día 30 de julio del año en curso en terminos del art. 58 de la Ley General de Contabilidad Gubernamental."

# Row 3 header area grew slightly taller
$ws.Rows("3:3").RowHeight = 40.5
# Row 8 data row height
$ws.Rows("8:8").RowHeight = 49.5

# Column width adjustments (F/G shrank now that long hyperlinks are gone, K grew for the note)
$ws.Columns("F:F").ColumnWidth = 65.140625
$ws.Columns("G:G").ColumnWidth = 80.5703125
$ws.Columns("K:K").ColumnWidth = 80.7109375

# Shrink the data-validation list range on column D
$ws.Range("D8:D97").Validation.Delete()
$ws.Range("D8:D97").Validation.Add(3, 1, 1, "=Hidden_13")

# View tidy-up
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("A14").Select()

$ws2.Range("A1").Select()
